$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Update Priority column (E) for rows 4-7 from "low" to "ht" on both locale sheets
foreach ($r in 4..7) {
    $zhcn.Range("E$r").Value = "ht"
    $dede.Range("E$r").Value = "ht"
}

# Update "Latest Handoff Datetime" column (H) for rows 4-7
foreach ($r in 4..7) {
    $zhcn.Range("H$r").Value = "2016-09-06 15:03:48"
    $dede.Range("H$r").Value = "2016-09-06 15:04:00"
}

# Update "Latest HO Xliff Generate Date" column (G) on Overview sheet for rows 4-7
foreach ($r in 4..7) {
    $overview.Range("G$r").Value = "2016-09-06 15:04:00"
}
